$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHIVE")

# Row 19: status moves to "In Progress" and hours logged
$ws.Range("E19").Value = "In Progress"
$ws.Range("G19").Value = 1.5

# Row 27: status moves to "In Progress" and hours revised down
$ws.Range("E27").Value = "In Progress"
$ws.Range("G27").Value = 1.5

# Row 32: "Implement deterministic world generation using seeds" removed from
# backlog, so the row now holds the task that used to be on row 33 and is
# marked complete
$ws.Range("B32").Value = "Task"
$ws.Range("D32").Value = "Write tests for world generation consistency"
$ws.Range("E32").Value = "Complete"
$ws.Range("G32").Value = 1

# Row 33: shifts up to the former row-34 task, now in progress
$ws.Range("B33").Value = "Feature"
$ws.Range("D33").Value = "Implement player movement and interaction"
$ws.Range("E33").Value = "In Progress"
$ws.Range("G33").Value = 2

# Row 34: shifts up to the former row-35 task, now in progress
$ws.Range("D34").Value = "Implement inventory system and item management"
$ws.Range("E34").Value = "In Progress"
$ws.Range("G34").Value = 3

# Row 35: new sprint 3 item - enemy movement and combat system
$ws.Range("D35").Value = "Implement Enemy movement and combact system"
$ws.Range("E35").Value = "In Progress"
$ws.Range("G35").Value = 3

# Row 37: "Implement combat or challenge mechanics (if part of the design)"
# removed from backlog and replaced by the trimmed version, now in progress
$ws.Range("D37").Value = "Implement combat or challenge mechanics "
$ws.Range("E37").Value = "In Progress"
$ws.Range("G37").Value = 2

# Row 38: hours logged
$ws.Range("G38").Value = 1

$ws.Range("D38").Select()
